# "Generate Report for Handback"
#
# This script reproduces the handback-report generation step for the
# aedb0925-6683-4b18-969d-b430955bf0a3 work item: for both the "zh-cn"
# and "de-de" localization sheets, it records that a handback file was
# received (row 8), links to the corresponding target xliff file, notes
# that the handback is stale, and records the error detail plus the
# datetime the check was performed. It also widens the "Error Detail"
# column so the long message is readable.

$wb = $excel.ActiveWorkbook

$targetMdAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7936d4d4299ee090adc278af8179a81530bbc2df/e2e/aedb0925-6683-4b18-969d-b430955bf0a3.md"
$targetMdDisplay = "aedb0925-6683-4b18-969d-b430955bf0a3.md"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/82def30dcb50135f87fbd30aa7b884ae9f357f7e/e2e/aedb0925-6683-4b18-969d-b430955bf0a3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7936d4d4299ee090adc278af8179a81530bbc2df/e2e/aedb0925-6683-4b18-969d-b430955bf0a3.md."

# Hyperlink font color used by the workbook's built-in "HyperLink" cell
# style (RGB 0x6495ED, stored as an OLE (BGR) color value for COM).
$hyperlinkColor = 15570276

function Set-HandbackRow8 {
    param(
        [string]$SheetName,
        [string]$TargetFileValue,
        [string]$HandbackDateTimeValue
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Widen column P ("Error Detail") so the long message is visible.
    $ws.Columns.Item(16).ColumnWidth = 39.1666666666667

    # I8 - "Latest Target File": new hyperlink to the handback markdown file.
    $ws.Range("I8").Value = $targetMdDisplay
    $ws.Hyperlinks.Add($ws.Range("I8"), $targetMdAddress, [Type]::Missing, [Type]::Missing, $targetMdDisplay) | Out-Null
    $ws.Range("I8").Font.Name = "Calibri"
    $ws.Range("I8").Font.Size = 11
    $ws.Range("I8").Font.Underline = 2
    $ws.Range("I8").Font.Color = $hyperlinkColor

    # J8 - "Latest Handback File": the generated target xliff file name.
    $ws.Range("J8").Value = $TargetFileValue

    # K8 - "Latest Handback DateTime": when the handback was processed.
    $ws.Range("K8").Value = $HandbackDateTimeValue

    # P8 - "Error Detail": the staleness warning for this handback.
    $ws.Range("P8").Value = $errorDetail
}

Set-HandbackRow8 -SheetName "zh-cn" -TargetFileValue "aedb0925-6683-4b18-969d-b430955bf0a3.4cdf1b1f80a97610f46c9b6535873e2ef7008085.zh-cn.xlf" -HandbackDateTimeValue "2016-08-30 02:45:35"

Set-HandbackRow8 -SheetName "de-de" -TargetFileValue "aedb0925-6683-4b18-969d-b430955bf0a3.4cdf1b1f80a97610f46c9b6535873e2ef7008085.de-de.xlf" -HandbackDateTimeValue "2016-08-30 02:45:42"
